$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 120
$ws.Range("I2").Value = 355
$ws.Range("J2").Value = 1505
$ws.Range("K2").Value = 12
$ws.Range("L2").Value = 395
$ws.Range("M2").Value = 21
$ws.Range("N2").Value = 271
$ws.Range("Q2").Value = 6
$ws.Range("R2").Value = 21
$ws.Range("S2").Value = 178
$ws.Range("T2").Value = 263
$ws.Range("V2").Value = 2440
$ws.Range("X2").Value = 2357
$ws.Range("Y2").Value = 1
$ws.Range("AA2").Value = 16
